$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1 title.
#    We clone the formatted content of the existing bold-title paragraph
#    (near the end of the doc) so that the resulting run layout matches
#    (leading empty run + bold run), then restyle/retext it.
# ------------------------------------------------------------------

$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$count = $d.Paragraphs.Count
$boldTitlePara = $d.Paragraphs.Item($count - 1)
$clonedFormattedText = $boldTitlePara.Range.FormattedText

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.FormattedText = $clonedFormattedText

$metaPara = $d.Paragraphs.Item(2)
$metaRange = $metaPara.Range
$metaRange.Find.Execute("Play Book of the Sphinx for Free " + [char]8211 + " Review & Bonus Offers", $true, $false, $false, $false, $false, $true, 1, $false, "Meta description", 2)

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.InsertAfter(": Read our review of Book of the Sphinx online slot game. Play for free and discover bonus offers. Suitable for all players. Try your luck now!")

# ------------------------------------------------------------------
# 2) Remove the duplicated bold title paragraph that used to precede the
#    italic meta-description paragraph at the end of the document.
# ------------------------------------------------------------------

$count = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs.Item($count - 1)
$dupTitlePara.Range.Delete()

# ------------------------------------------------------------------
# 3) Replace the text of the trailing italic paragraph with the new
#    image-generation prompt, keeping its italic run formatting.
#    (Use Find to locate + Range.Text = ... instead of Find.Execute's
#    replace argument, which smart-quotes curly apostrophes.)
# ------------------------------------------------------------------

$count = $d.Paragraphs.Count
$promptPara = $d.Paragraphs.Item($count)
$promptRange = $promptPara.Range
$promptRange.Find.Execute("Read our review of Book of the Sphinx online slot game. Play for free and discover bonus offers. Suitable for all players. Try your luck now!")

$newPromptText = 'Please create a cartoon-style feature image for "Book of the Sphinx" with a happy Maya warrior with glasses. The warrior should be standing in front of the pyramids, holding the Book of the Sphinx in one hand and waving the other hand in excitement. He should be wearing a traditional Maya warrior outfit, complete with a headdress and a weapon. The background should feature the sunset over the pyramids, with warm shades of orange and yellow. The overall tone of the image should be joyful and playful, capturing the fun and adventurous spirit of the game. Please ensure that the image is vibrant, eye-catching, and in line with the game''s theme.'
$promptRange.Text = $newPromptText

Write-Output "done"
